# Update the division-practice answer table: each data cell's text is
# replaced with its new value. Every "old" string below is unique across
# the whole document, so a plain Find/Execute (wdReplaceAll) per pair is
# safe and unambiguous, even though a couple of "new" values happen to
# equal an "old" value used earlier in the list (those earlier entries
# are already consumed by the time the matching value is produced).
$d = $word.ActiveDocument

$d.Content.Find.Execute("717÷5=143, 2", $true, $false, $false, $false, $false, $true, 1, $false, "814÷6=135, 4", 2)
$d.Content.Find.Execute("146÷4=36, 2", $true, $false, $false, $false, $false, $true, 1, $false, "411÷9=45, 6", 2)
$d.Content.Find.Execute("945÷6=157, 3", $true, $false, $false, $false, $false, $true, 1, $false, "470÷6=78, 2", 2)
$d.Content.Find.Execute("964÷7=137, 5", $true, $false, $false, $false, $false, $true, 1, $false, "542÷6=90, 2", 2)
$d.Content.Find.Execute("493÷9=54, 7", $true, $false, $false, $false, $false, $true, 1, $false, "302÷5=60, 2", 2)
$d.Content.Find.Execute("779÷8=97, 3", $true, $false, $false, $false, $false, $true, 1, $false, "576÷6=96, 0", 2)
$d.Content.Find.Execute("109÷3=36, 1", $true, $false, $false, $false, $false, $true, 1, $false, "130÷3=43, 1", 2)
$d.Content.Find.Execute("647÷6=107, 5", $true, $false, $false, $false, $false, $true, 1, $false, "459÷7=65, 4", 2)
$d.Content.Find.Execute("634÷5=126, 4", $true, $false, $false, $false, $false, $true, 1, $false, "125÷4=31, 1", 2)
$d.Content.Find.Execute("884÷2=442, 0", $true, $false, $false, $false, $false, $true, 1, $false, "782÷8=97, 6", 2)
$d.Content.Find.Execute("568÷7=81, 1", $true, $false, $false, $false, $false, $true, 1, $false, "493÷9=54, 7", 2)
$d.Content.Find.Execute("212÷4=53, 0", $true, $false, $false, $false, $false, $true, 1, $false, "895÷5=179, 0", 2)
$d.Content.Find.Execute("317÷4=79, 1", $true, $false, $false, $false, $false, $true, 1, $false, "162÷2=81, 0", 2)
$d.Content.Find.Execute("744÷6=124, 0", $true, $false, $false, $false, $false, $true, 1, $false, "703÷6=117, 1", 2)
$d.Content.Find.Execute("562÷6=93, 4", $true, $false, $false, $false, $false, $true, 1, $false, "833÷7=119, 0", 2)
$d.Content.Find.Execute("928÷8=116, 0", $true, $false, $false, $false, $false, $true, 1, $false, "541÷8=67, 5", 2)
$d.Content.Find.Execute("923÷2=461, 1", $true, $false, $false, $false, $false, $true, 1, $false, "627÷3=209, 0", 2)
$d.Content.Find.Execute("192÷8=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "992÷2=496, 0", 2)
$d.Content.Find.Execute("693÷2=346, 1", $true, $false, $false, $false, $false, $true, 1, $false, "845÷3=281, 2", 2)
$d.Content.Find.Execute("930÷5=186, 0", $true, $false, $false, $false, $false, $true, 1, $false, "562÷6=93, 4", 2)
$d.Content.Find.Execute("781÷9=86, 7", $true, $false, $false, $false, $false, $true, 1, $false, "859÷5=171, 4", 2)
$d.Content.Find.Execute("202÷6=33, 4", $true, $false, $false, $false, $false, $true, 1, $false, "710÷4=177, 2", 2)
$d.Content.Find.Execute("954÷8=119, 2", $true, $false, $false, $false, $false, $true, 1, $false, "320÷6=53, 2", 2)
$d.Content.Find.Execute("458÷8=57, 2", $true, $false, $false, $false, $false, $true, 1, $false, "644÷2=322, 0", 2)
$d.Content.Find.Execute("809÷8=101, 1", $true, $false, $false, $false, $false, $true, 1, $false, "940÷2=470, 0", 2)
